$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'report'
$ws.Range("A2").Value = 'Grade 10 student believed a black van followed as she walked to school from her resident on three occasions. This morning as the victim was walking to school. An unknown suspect jumped out from the bush and groped victim from behind, then ran away. The parents were contacted and a safety plan has been put in place including for victim to not walk to school alone in the mornings for the time being. Preliminary area search negative for the vehicle and several CCTV potential locations observed.'
$ws.Range("A3").Value = 'Victim was walking home. An unknown SOC male approached her from behind, groped her buttocks, and ran away westbound on Pandora Ave. Victim called police. Victim stated this was the second time this had occurred in two weeks with the previous incident not being reported. Police attended and conducted an area search with negative results. Video canvas conducted with possible sources located. Safety plan created, and victim services offered.'
$ws.Range("A4").Value = 'Victim was waiting at a bus stop in the 3500 blk Kingsway ave when an unknown male approached her from behind and groped her buttocks. The male then followed victim in to the 3300 blk Padora Ave and fled the area. Written statement and video of suspect obtained. '
$ws.Range("A5").Value = 'Police responded to a sex assault that occurred in the west lane of  200 Pandora St.  Victim was groped from behind while she waited at the rear door of her office building.  The male suspect approached her from behind and grabbed her buttox/groin area with one hand while wrapping his other arm around her like a hug/embrace.  While the male suspect walked away he looked back at the victim and smiled.  In response, victim ran after the suspect and filmed him with her cell phone.  Video canvas of the area captured the suspect and the assault on CCTV. '
$ws.Range("A6").Value = 'Victim was standing in front of Fresh Donair located at 312 Main st in Downtown. Vancouver when an unknown South Asian male groped her butt from behind. Victim started to video record the male and attempted to confront him as he pulled up his hoody and started to run away. The male fled NB toward Nelson st. At approx. 1343 hrs, police were dispatched to the above file and arranged an audio recorded interview with victim. PCs obtained a statement from victim and seized her long coat for evidence. Victim was offered victim services and video canvas to be conducted. Victim believed suspect was hiding in the bush for a while.'
$ws.Range("A7").Value = 'Police attempted to conduct a traffic stop with a black 2010 Toyota Corolla bearing BCLP: ABC123 in the 1200 bl E.41st Ave for driving in the bus lane. Emergency equipment was activated and the vehicle which had 2 Asian occupants made a abrupt right turn onto Sherbrooke St, striking the curb on the east side of the street and came to a stop. '
$ws.Range("A8").Value = 'Vehicle 1 a blue 2017 toyota rav4 driven by SIMPSON was driving eastbound on Nelson Street approaching Granville Street. Vehicle 2 A 2017 white Ford F150 registered to LOL was also driving eastbound on Nelson Street. In the intersection at Granville Street and Nelson Street the two vehicle were side by side when they collided (vehicle 1''s front passenger side hit the front drivers side of Vehicle 2). Both vehicles pulled over briefly before Vehicle 2 driven by a caucasian male with a female in the car drove off without exchanging information.'
$ws.Range("A9").Value = 'Witness heard the loud engine of his truck rev up and speed off quickly. Vehicle was last seen travelling Northbound on Nanaimo St. Witness was notified his truck was stolen. witness attempted to locate his vehicle, however was unsuccessful. 1986 Yellow Chevy pick up truck, with wooden panels on the side (Junk removal truck). Trunk door has "ladies and gentleman" painted in the back. '
$ws.Range("A10").Value = 'Suspects robbed a bank and got away with vehicle. Vehicle description Acura SUV white, older model.'
$ws.Range("A11").Value = 'Victim Grade 11 at Killarney Secondary School was walking to a grocery store at Champlain Square near the intersection of Pandora Ave/Kerr Street, Vancouver, BC, when all of a sudden her crotch was grabbed over her shorts from behind by an unknown male, who then fled the area on foot. On April 21, 2021 after she finished her afternoon class, victim approached her SLO in the school hallway and advised him of the circumstances. Victim did not observe the male''s face, but observed him to be wearing a black hoody with the hood up and black or grey sweatpants. Video canvass to be conducted and victim''s unwashed shorts to be obtained and submitted for forensic processing. A safety plan was implemented and support was provided by her school counsellor in addition to the request for Victim Services at a later date.'

$ws.Range("A3").Select()
